$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 10
$ws.Range("E3").Value = 20

$ws.Range("E3").Select()
